# Add two new columns (I: "I0", J: "IF") to Sheet1.
# I is a constant 1 for every data row, J mirrors the existing H column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells, with the same style/formatting as the existing header row
# (bold font, centered/top aligned, bordered) copied from H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows: I = 1 (constant), J = same value as H for that row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
